# Apply the template change: remove the "11 DIGITS" column (column C) from the
# MPD sheet. Deleting the entire column shifts the subsequent columns left,
# which matches the target diff (old col D -> new col C, etc.) and also drops
# the now-unused "11 DIGITS" shared string automatically.

$wb = $excel.ActiveWorkbook

$mpd = $wb.Worksheets.Item("MPD")

# Delete the entire column C ("11 DIGITS") - shifts remaining columns left.
$mpd.Columns.Item(3).Delete()

# The hidden "_FilterDatabase" name for MPD tracked the old A1:L1 header
# range; bring it in line with the now-narrower A1:K1 header row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "MPD!_FilterDatabase") {
        $n.RefersTo = "=MPD!`$A`$1:`$K`$1"
    }
}

$mpd.Activate()
$mpd.Range("F6").Select()
